$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "https://www.myntra.com/bra/amante/amante-solid-padded-wired-multiway-strapless-bra---bra10808/2528934/buy"
$ws.Range("B2").Value = "Amante"
$ws.Range("C2").Value = "Solid Padded Wired Multiway Strapless Bra - BRA10808"
$ws.Range("D2").Value = "'4.2"
$ws.Range("E2").Value = "663 Ratings"

# Row 3
$ws.Range("A3").Value = "https://www.myntra.com/bra/marks+%26+spencer/marks--spencer-bra-full-coverage-underwired-lightly-padded/32878492/buy"
$ws.Range("B3").Value = "Marks & Spencer"
$ws.Range("C3").Value = "Bra Full Coverage Underwired Lightly Padded"
$ws.Range("D3").Value = "'4.5"
$ws.Range("E3").Value = "33 Ratings"

# Row 4
$ws.Range("A4").Value = "https://www.myntra.com/bra/amante/amante-floral-bra-full-coverage/28111984/buy"
$ws.Range("B4").Value = "Amante"
$ws.Range("C4").Value = "Floral Bra Full Coverage"
$ws.Range("D4").Value = "'4.3"
$ws.Range("E4").Value = "357 Ratings"

# Row 5
$ws.Range("A5").Value = "https://www.myntra.com/bra/amante/amante-floral-bra-full-coverage/28111990/buy"
$ws.Range("B5").Value = "Amante"
$ws.Range("C5").Value = "Floral Bra Full Coverage"
$ws.Range("E5").Value = "623 Ratings"

# Delete rows 6-12
$ws.Range("A6:E12").EntireRow.Delete()
